$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update formula in N2 (8*15360 -> 7.5*15360)
$ws.Range("N2").Formula = "=7.5*15360"

# Update the active cell selection to N2
$ws.Range("N2").Select()
